$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6830.208516526744
$ws.Range("C2").Value = 6564.051123591939
$ws.Range("D2").Value = 26449.14472140368
